$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 8
$ws.Range("C3").Value = 8

$ws.Range("C4").Value = 8

$ws.Range("C5").Value = 12

$ws.Range("A6").Value = 15
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 17

$ws.Range("A7").Value = 20
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 16

$ws.Range("B8").Select()
